$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows above the existing data (new rows 2-9); existing rows 2-21 shift down to 10-29.
$ws.Rows("2:9").Insert()
$ws.Rows("2:9").ClearFormats()

# Values for the newly inserted rows 2-9 (x, y, z)
$newTopData = @(
    @(-0.0337503030896186, 0.09865473955869671, -0.0665843114256858),
    @(-0.0201585534960031, -0.0109955742955207, -0.0201585534960031),
    @(0.0261144898831844, -0.0386372283101081, 0.0106901414692401),
    @(0.0059559359215199, -0.0545197241008281, 0.0339030213654041),
    @(0.0154243474826216, -0.0274889357388019, -0.0103847095742821),
    @(-0.0024434609804302, -0.0116064399480819, -0.0108428578823804),
    @(0.0427605658769607, 0.0514653958380222, -0.0694859251379966),
    @(0.0119118718430399, -0.0441350154578685, 0.0128281703218817),
)

$r = 2
foreach ($row in $newTopData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Append 2 new rows (30-31) after the existing data
$newBottomData = @(
    @(0.0019853119738399, -0.0036651915870606, 0.0067195175215601),
    @(0.0003054326225537, -0.00167987938039, 0.0142026171088218),
)

$r = 30
foreach ($row in $newBottomData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
